$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 112038436
$ws.Range("AB2").Value2 = "13:28"
$ws.Range("AM2").Value2 = "Stående död trädstam/högstubbe"
$ws.Range("AO2").Value2 = "Standing dead tree/snags # Picea abies"
$ws.Range("B2").Value2 = 89535
$ws.Range("E2").Value2 = 1108
$ws.Range("F2").Value2 = "Harticka"
$ws.Range("G2").Value2 = "Pelloporus leporinus"
$ws.Range("H2").Value2 = "(Fr.) Krieglst."
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("Q2").Value2 = 515951
$ws.Range("Z2").Value2 = "13:28"

# Row 3
$ws.Range("A3").Value2 = 112038473
$ws.Range("AB3").Value2 = "13:34"
$ws.Range("AM3").Value2 = "Liggande död trädstam, markontakt"
$ws.Range("AO3").Value2 = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("B3").Value2 = 89820
$ws.Range("E3").Value2 = 658
$ws.Range("F3").Value2 = "Rosenticka"
$ws.Range("G3").Value2 = "Rhodofomes roseus"
$ws.Range("H3").Value2 = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value2 = "4"
$ws.Range("J3").Value2 = "fruktkroppar"
$ws.Range("Q3").Value2 = 516057
$ws.Range("Z3").Value2 = "13:34"

# Row 7
$ws.Range("A7").Value2 = 112035020
$ws.Range("AB7").Value2 = "10:24"
$ws.Range("AJ7").ClearContents()
$ws.Range("AK7").ClearContents()
$ws.Range("AM7").Value2 = "Stubbe"
$ws.Range("AO7").Value2 = "Stump"
$ws.Range("B7").Value2 = 89535
$ws.Range("E7").Value2 = 1108
$ws.Range("F7").Value2 = "Harticka"
$ws.Range("G7").Value2 = "Pelloporus leporinus"
$ws.Range("H7").Value2 = "(Fr.) Krieglst."
$ws.Range("Q7").Value2 = 515923
$ws.Range("R7").Value2 = 7184659
$ws.Range("S7").Value2 = 50
$ws.Range("Z7").Value2 = "10:24"

# Row 8
$ws.Range("A8").Value2 = 112037386
$ws.Range("AB8").Value2 = "11:52"
$ws.Range("AH8").Value2 = "Blåbärsbarrskog"
$ws.Range("B8").Value2 = 89557
$ws.Range("E8").Value2 = 5432
$ws.Range("F8").Value2 = "Granticka"
$ws.Range("G8").Value2 = "Porodaedalea chrysoloma"
$ws.Range("H8").Value2 = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q8").Value2 = 516032
$ws.Range("R8").Value2 = 7184227
$ws.Range("Z8").Value2 = "11:52"

# Row 9
$ws.Range("A9").Value2 = 112037635
$ws.Range("AB9").Value2 = "12:06"
$ws.Range("AH9").Value2 = "Blåbärsgranskog"
$ws.Range("AJ9").Value2 = "gran"
$ws.Range("AK9").Value2 = "Picea abies"
$ws.Range("AM9").Value2 = "Liggande död trädstam, markontakt"
$ws.Range("AO9").Value2 = "Horizontal, dead with ground contact # Picea abies"
$ws.Range("Q9").Value2 = 515886
$ws.Range("R9").Value2 = 7184226
$ws.Range("S9").Value2 = 10
$ws.Range("Z9").Value2 = "12:06"
